# Update "想去人数" (interest count) figures in the 苏州-漫展信息 workbook.
# Mirrors the upstream gh-pages data refresh commit 456a3b4: a handful of
# event rows on the 展览 / 演出 / 全部类型 sheets had their F-column counts
# bumped up by the latest scrape.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 111
$ws1.Range("F7").Value  = 11813
$ws1.Range("F8").Value  = 4417
$ws1.Range("F14").Value = 1103
$ws1.Range("F15").Value = 158
$ws1.Range("F16").Value = 47
$ws1.Range("F17").Value = 5131
$ws1.Range("F19").Value = 189
$ws1.Range("F20").Value = 529
$ws1.Range("F21").Value = 11371
$ws1.Range("F22").Value = 11334
$ws1.Range("F23").Value = 19
$ws1.Range("F24").Value = 50
$ws1.Range("F25").Value = 13
$ws1.Range("F27").Value = 49

# --- 演出 sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# --- 全部类型 sheet ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 111
$ws4.Range("F7").Value  = 11813
$ws4.Range("F8").Value  = 4417
$ws4.Range("F14").Value = 2
$ws4.Range("F15").Value = 1103
$ws4.Range("F16").Value = 158
$ws4.Range("F17").Value = 47
$ws4.Range("F18").Value = 5131
$ws4.Range("F20").Value = 189
$ws4.Range("F21").Value = 529
$ws4.Range("F22").Value = 11371
$ws4.Range("F23").Value = 11334
$ws4.Range("F24").Value = 19
$ws4.Range("F25").Value = 50
$ws4.Range("F26").Value = 13
$ws4.Range("F28").Value = 49
